$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "2026-02-18 14:45:07"
    "2026-02-18 12:17:29"
    "2026-02-18 15:03:18"
    "2026-02-18 18:15:30"
    "2026-02-18 12:47:28"
    "2026-02-18 13:49:45"
    "2026-02-17 09:06:24"
    "2026-02-18 17:14:02"
    "2026-02-18 11:04:26"
    "2026-02-18 12:47:22"
    "2026-02-18 10:48:06"
)
$numbers = @(
    "237679041654"
    "237673671238"
    "237652275301"
    "237681662701"
    "237654349065"
    "237675637054"
    "237671262234"
    "237674890488"
    "237654079053"
    "237652643069"
    "237673041651"
)
$names = @(
    "PIERRE MARIVOT TEMEZEU"
    "LA NEGRESSE SARL MBONE NDEMOU EPSE KAMSU ROSINE"
    "NDAMI EPSE NONGA ROSALIE ETS MOBILE FINANCIAL SERVICES MFS"
    "TOUGOUA PAYOU JULIO OMER ETS MOBILE FINANCIAL SERVICES MFS"
    "YASSI A BAA BELMOND CHIC MOBILE"
    "N A SOKOUDJOU DZOKOU"
    "DEUGOUE TOKO EPSE DIBANGUE LOISE LAURE ETS MOBILE FINANCIAL SERVICES MFS"
    "Sandrine Nkendji"
    "JUDITH AIMEE JOELEFACK JAZET EPSE NGUMATIO"
    "CHIREL DELRICH TCHAPDA"
    "DYLAN KEPSEU SIME"
)
$balances = @(6876, 4750, 7688, 17160, 360662, 6316, 70, 253822, 216137, 6503, 817810)

$startRow = 178
for ($i = 0; $i -lt $dates.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $dates[$i]
}
for ($i = 0; $i -lt $numbers.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = "'" + $numbers[$i]
    $ws.Cells.Item($startRow + $i, 2).Style = "Normal"
}
for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $names[$i]
}
for ($i = 0; $i -lt $balances.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $balances[$i]
}
